$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.069404619044566
$ws.Range("D2").Value = 1.074188088067272
$ws.Range("E2").Value = 1.063392055330244
$ws.Range("F2").Value = 1.082450564687025
$ws.Range("J2").Value = 1.074339515050226
$ws.Range("K2").Value = 1.076877962343133
$ws.Range("L2").Value = 1.066110819384371
$ws.Range("M2").Value = 1.085118749329561
$ws.Range("N2").Value = 1.075865200390062
$ws.Range("C3").Value = 1.073499821534698
$ws.Range("D3").Value = 1.078068119455432
$ws.Range("E3").Value = 1.067014910477174
$ws.Range("F3").Value = 1.086553944373795
$ws.Range("J3").Value = 1.078077137737485
$ws.Range("K3").Value = 1.08056757023439
$ws.Range("L3").Value = 1.069541647556044
$ws.Range("M3").Value = 1.089032851632673
$ws.Range("N3").Value = 1.079608130930248
$ws.Range("C4").Value = 1.076122029012796
$ws.Range("D4").Value = 1.080551921490236
$ws.Range("E4").Value = 1.069333299430811
$ws.Range("F4").Value = 1.089181731879241
$ws.Range("J4").Value = 1.080468739226109
$ws.Range("K4").Value = 1.082928205225848
$ws.Range("L4").Value = 1.071735790851134
$ws.Range("M4").Value = 1.091538249183221
$ws.Range("N4").Value = 1.082003128767306
$ws.Range("C5").Value = 1.077218013818547
$ws.Range("D5").Value = 1.081589905841725
$ws.Range("E5").Value = 1.070301969291542
$ws.Range("F5").Value = 1.090280129657587
$ws.Range("J5").Value = 1.081467945314063
$ws.Range("K5").Value = 1.083914415175479
$ws.Range("L5").Value = 1.072652223395527
$ws.Range("M5").Value = 1.092585207443027
$ws.Range("N5").Value = 1.083003753842516
$ws.Range("C6").Value = 1.077401666291574
$ws.Range("D6").Value = 1.081763830294938
$ws.Range("E6").Value = 1.070464268354773
$ws.Range("F6").Value = 1.090464191158566
$ws.Range("J6").Value = 1.081635357590057
$ws.Range("K6").Value = 1.084079646577313
$ws.Range("L6").Value = 1.072805751137326
$ws.Range("M6").Value = 1.092760632593145
$ws.Range("N6").Value = 1.083171403863145
$ws.Range("C7").Value = 1.076136698426808
$ws.Range("D7").Value = 1.080565815187576
$ws.Range("E7").Value = 1.069346266077436
$ws.Range("F7").Value = 1.08919643327435
$ws.Range("J7").Value = 1.080482114839759
$ws.Range("K7").Value = 1.082941407100331
$ws.Range("L7").Value = 1.071748059526135
$ws.Range("M7").Value = 1.091552263198381
$ws.Range("N7").Value = 1.082016523375862
$ws.Range("C8").Value = 1.070794481193901
$ws.Range("D8").Value = 1.075505053895055
$ws.Range("E8").Value = 1.064621892479433
$ws.Range("F8").Value = 1.0838431293129
$ws.Range("J8").Value = 1.075608361290322
$ws.Range("K8").Value = 1.07813055916407
$ws.Range("L8").Value = 1.067275754349049
$ws.Range("M8").Value = 1.086447324754199
$ws.Range("N8").Value = 1.077135848537358
$ws.Range("C9").Value = 1.061158198065161
$ws.Range("D9").Value = 1.066371616934374
$ws.Range("E9").Value = 1.056089480605242
$ws.Range("F9").Value = 1.074189616686943
$ws.Range("J9").Value = 1.066804264487162
$ws.Range("K9").Value = 1.069438200897037
$ws.Range("L9").Value = 1.059187887108215
$ws.Range("M9").Value = 1.07723244901418
$ws.Range("N9").Value = 1.068319248906892
$ws.Range("C10").Value = 1.054568734874692
$ws.Range("D10").Value = 1.060122839671849
$ws.Range("E10").Value = 1.050247859037545
$ws.Range("F10").Value = 1.067590356509285
$ws.Range("J10").Value = 1.060775224264203
$ws.Range("K10").Value = 1.063484407765978
$ws.Range("L10").Value = 1.053643312307956
$ws.Range("M10").Value = 1.070926757649121
$ws.Range("N10").Value = 1.062281646755276
$ws.Range("C11").Value = 1.051672573472101
$ws.Range("D11").Value = 1.057375686463261
$ws.Range("E11").Value = 1.047678757758961
$ws.Range("F11").Value = 1.064690386353461
$ws.Range("J11").Value = 1.058123332759835
$ws.Range("K11").Value = 1.060865314548612
$ws.Range("L11").Value = 1.05120309939799
$ws.Range("M11").Value = 1.068154289907685
$ws.Range("N11").Value = 1.059625989260796
$ws.Range("C12").Value = 1.050590041262798
$ws.Range("D12").Value = 1.056348743558179
$ws.Range("E12").Value = 1.046718233444627
$ws.Range("F12").Value = 1.063606507476275
$ws.Range("J12").Value = 1.057131797778756
$ws.Range("K12").Value = 1.059885997993526
$ws.Range("L12").Value = 1.050290498862824
$ws.Range("M12").Value = 1.067117839014086
$ws.Range("N12").Value = 1.058633046186314
$ws.Range("C13").Value = 1.050822559818209
$ws.Range("D13").Value = 1.056569326824624
$ws.Range("E13").Value = 1.046924556611527
$ws.Range("F13").Value = 1.063839311766787
$ws.Range("J13").Value = 1.057344784757568
$ws.Range("K13").Value = 1.060096362409732
$ws.Range("L13").Value = 1.050486539853503
$ws.Range("M13").Value = 1.067340466601679
$ws.Range("N13").Value = 1.058846335631066
$ws.Range("C14").Value = 1.051583231047804
$ws.Range("D14").Value = 1.057290934048185
$ws.Range("E14").Value = 1.047599489655193
$ws.Range("F14").Value = 1.064600931206836
$ws.Range("J14").Value = 1.058041506658048
$ws.Range("K14").Value = 1.060784497678436
$ws.Range("L14").Value = 1.051127791637374
$ws.Range("M14").Value = 1.068068753727957
$ws.Range("N14").Value = 1.05954404695656
$ws.Range("C15").Value = 1.052050998725617
$ws.Range("D15").Value = 1.057734665573341
$ws.Range("E15").Value = 1.048014501656273
$ws.Range("F15").Value = 1.065069292219917
$ws.Range("J15").Value = 1.058469908820806
$ws.Range("K15").Value = 1.061207614162276
$ws.Range("L15").Value = 1.05152205826036
$ws.Range("M15").Value = 1.068516586909231
$ws.Range("N15").Value = 1.059973057499528
$ws.Range("C16").Value = 1.054760012593406
$ws.Range("D16").Value = 1.060304261038091
$ws.Range("E16").Value = 1.05041750218432
$ws.Range("F16").Value = 1.067781896382617
$ws.Range("J16").Value = 1.060950326364016
$ws.Range("K16").Value = 1.063657337972954
$ws.Range("L16").Value = 1.053804407909191
$ws.Range("M16").Value = 1.07110984456905
$ws.Range("N16").Value = 1.062456997520156
$ws.Range("C17").Value = 1.056447599266
$ws.Range("D17").Value = 1.061904804836991
$ws.Range("E17").Value = 1.051914027355046
$ws.Range("F17").Value = 1.06947185322566
$ws.Range("J17").Value = 1.062494966149789
$ws.Range("K17").Value = 1.065182784420728
$ws.Range("L17").Value = 1.055225329233717
$ws.Range("M17").Value = 1.072725048573062
$ws.Range("N17").Value = 1.064003830871598
$ws.Range("C18").Value = 1.057427825519445
$ws.Range("D18").Value = 1.06283440290319
$ws.Range("E18").Value = 1.052783120706712
$ws.Range("F18").Value = 1.070453503568833
$ws.Range("J18").Value = 1.063391966919018
$ws.Range("K18").Value = 1.066068610634144
$ws.Range("L18").Value = 1.056050349092324
$ws.Range("M18").Value = 1.073663133312905
$ws.Range("N18").Value = 1.064902105484807
$ws.Range("C19").Value = 1.05776136911943
$ws.Range("D19").Value = 1.063150707248724
$ws.Range("E19").Value = 1.05307882236829
$ws.Range("F19").Value = 1.070787539717096
$ws.Range("J19").Value = 1.063697158033274
$ws.Range("K19").Value = 1.066369994910177
$ws.Range("L19").Value = 1.05633102679533
$ws.Range("M19").Value = 1.073982320603663
$ws.Range("N19").Value = 1.06520773000545
$ws.Range("C20").Value = 1.056266965082055
$ws.Range("D20").Value = 1.06173349467292
$ws.Range("E20").Value = 1.051753859935153
$ws.Range("F20").Value = 1.069290960420143
$ws.Range("J20").Value = 1.062329652819211
$ws.Range("K20").Value = 1.065019528200514
$ws.Range("L20").Value = 1.055073270745956
$ws.Range("M20").Value = 1.072552172197074
$ws.Range("N20").Value = 1.063838282777128
$ws.Range("C21").Value = 1.051359421885189
$ws.Range("D21").Value = 1.057078621464795
$ws.Range("E21").Value = 1.047400913506776
$ws.Range("F21").Value = 1.064376840916159
$ws.Range("J21").Value = 1.057836521481384
$ws.Range("K21").Value = 1.060582040038853
$ws.Range("L21").Value = 1.050939132355749
$ws.Range("M21").Value = 1.067854477015278
$ws.Range("N21").Value = 1.059338770677433
$ws.Range("C22").Value = 1.048234518871114
$ws.Range("D22").Value = 1.054113985804676
$ws.Range("E22").Value = 1.044627753924193
$ws.Range("F22").Value = 1.061248198549691
$ws.Range("J22").Value = 1.05497371858355
$ws.Range("K22").Value = 1.057754430387524
$ws.Range("L22").Value = 1.048303833781944
$ws.Range("M22").Value = 1.064862307429477
$ws.Range("N22").Value = 1.056471902271116
$ws.Range("C23").Value = 1.049894930753248
$ws.Range("D23").Value = 1.055689297867532
$ws.Range("E23").Value = 1.046101398555661
$ws.Range("F23").Value = 1.062910554219863
$ws.Range("J23").Value = 1.056495031783046
$ws.Range("K23").Value = 1.059257066119123
$ws.Range("L23").Value = 1.049704365278723
$ws.Range("M23").Value = 1.066452275074799
$ws.Range("N23").Value = 1.057995375909852
$ws.Range("C24").Value = 1.056348598586833
$ws.Range("D24").Value = 1.061810914615462
$ws.Range("E24").Value = 1.051826244434646
$ws.Range("F24").Value = 1.069372710656523
$ws.Range("J24").Value = 1.062404363018575
$ws.Range("K24").Value = 1.065093308816834
$ws.Range("L24").Value = 1.055141991091542
$ws.Range("M24").Value = 1.072630300041994
$ws.Range("N24").Value = 1.063913099073544
$ws.Range("C25").Value = 1.063677354202723
$ws.Range("D25").Value = 1.068759876566056
$ws.Range("E25").Value = 1.058321289480977
$ws.Range("F25").Value = 1.076712939882898
$ws.Range("J25").Value = 1.069107373762431
$ws.Range("K25").Value = 1.071712301888703
$ws.Range("L25").Value = 1.061304683107325
$ws.Range("M25").Value = 1.079642208938691
$ws.Range("N25").Value = 1.070625628861503
